# Weekly update: insert a new week's price observation as the new row 2,
# pushing all existing records down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits the row-above's (header) bold/centered
# style; strip it back to the default "no style" cells used by the rest of
# the data rows.
$ws.Range("A2:T2").ClearFormats()

# Restore the date number format on the Fecha column, matching the other
# data rows (style index shared with D3:D22).
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44922
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103001
$ws.Range("J2").Value = "Cereza"
$ws.Range("K2").Value = "Bing"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 5000
$ws.Range("O2").Value = 6000
$ws.Range("P2").Value = 5500
$ws.Range("Q2").Value = "`$/bandeja 10 kilos"
$ws.Range("R2").Value = "Región del Maule"
$ws.Range("S2").Value = 550
$ws.Range("T2").Value = 10
